$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values would otherwise be auto-parsed by Excel as
# numbers (losing the exact decimal text / becoming a float). Force them to
# stay plain text, matching the original inlineStr cell contents, then reset
# the style back to Normal so no stray formatting is introduced.
$textCells = @("D5", "D6", "D9", "D11", "D13", "D15", "D19", "D20", "D21", "D22", "D24", "D28", "D29", "D30", "D31", "D32", "D35", "D38", "D39", "D40", "D42", "D44", "D45", "D46", "D47", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price / volume figures.
$ws.Range("D2").Value = "69.430.37"
$ws.Range("E2").Value = "  -2.89%  "
$ws.Range("D3").Value = "3.686.45"
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "686.08"
$ws.Range("E5").Value = "  -2.38%  "
$ws.Range("D6").Value = "161.41"
$ws.Range("E6").Value = "  -6.01%  "
$ws.Range("D7").Value = "3.684.69"
$ws.Range("E7").Value = "  -3.51%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  -6.01%  "
$ws.Range("E10").Value = "  -8.75%  "
$ws.Range("D11").Value = "7.33"
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("E12").Value = "  -9.62%  "
$ws.Range("D13").Value = "0.0000235"
$ws.Range("E13").Value = "  -6.77%  "
$ws.Range("D14").Value = "4.305.63"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").Value = "32.90"
$ws.Range("E15").Value = "  -9.03%  "
$ws.Range("D16").Value = "3.685.97"
$ws.Range("E16").Value = "  -3.18%  "
$ws.Range("D17").Value = "69.457.01"
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").Value = "15.97"
$ws.Range("E19").Value = "  -9.18%  "
$ws.Range("D20").Value = "6.49"
$ws.Range("E20").Value = "  -10.25%  "
$ws.Range("D21").Value = "474.00"
$ws.Range("E21").Value = "  -8.01%  "
$ws.Range("D22").Value = "9.91"
$ws.Range("E22").Value = "  -5.54%  "
$ws.Range("E23").Value = "  -8.77%  "
$ws.Range("D24").Value = "79.75"
$ws.Range("E24").Value = "  -5.34%  "
$ws.Range("D25").Value = "3.831.83"
$ws.Range("E25").Value = "  -3.33%  "
$ws.Range("E26").Value = "  -10.03%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "11.20"
$ws.Range("E28").Value = "  -13.12%  "
$ws.Range("D29").Value = "9.20"
$ws.Range("E29").Value = "  -11.41%  "
$ws.Range("D30").Value = "1.78"
$ws.Range("E30").Value = "  -12.25%  "
$ws.Range("D31").Value = "2.69"
$ws.Range("E31").Value = "  -11.00%  "
$ws.Range("D32").Value = "6.74"
$ws.Range("E32").Value = "  -8.70%  "
$ws.Range("E33").Value = "  -9.23%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "26.74"
$ws.Range("E35").Value = "  -8.51%  "
$ws.Range("E36").Value = "  -5.39%  "
$ws.Range("D37").Value = "3.653.08"
$ws.Range("E37").Value = "  -3.38%  "
$ws.Range("D38").Value = "8.28"
$ws.Range("E38").Value = "  -10.41%  "
$ws.Range("D39").Value = "6.20"
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("D40").Value = "2.29"
$ws.Range("E40").Value = "  -7.32%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "0.0913"
$ws.Range("E42").Value = "  -9.98%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "0.949"
$ws.Range("E44").Value = "  -6.86%  "
$ws.Range("D45").Value = "163.91"
$ws.Range("E45").Value = "  -5.20%  "
$ws.Range("D46").Value = "48.28"
$ws.Range("E46").Value = "  -3.30%  "
$ws.Range("D47").Value = "30.09"
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("D48").Value = "2.72"
$ws.Range("E48").Value = "  -16.89%  "
$ws.Range("E49").Value = "  -4.58%  "
$ws.Range("D50").Value = "0.000278"
$ws.Range("E50").Value = "  -9.49%  "
$ws.Range("E51").Value = "  -5.54%  "

# Drop the temporary text format so styling matches the original workbook.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
